# Reorders the Item Name / UOM pairs within several brand groups on the
# active sheet (columns D and E), matching the updated source ordering.
# Only the rows whose Item Name/UOM pairing actually changes are touched;
# all other rows (and all other columns) are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @("Dinafex 180mg Tablet", "30's")
    5  = @("Dinafex 60mg Tablet", "30's")
    7  = @("Etorix 90mg Tablet", "30's")
    9  = @("Etorix 60mg Tablet - 40's", "40's")
    11 = @("Flucloxin 500mg Capsule", "30 's")
    12 = @("Flucloxin 500mg Capsule - 36's", "36 's")
    14 = @("Ketonic 30mg IM/IV Injection - 4's", "4's")
    15 = @("Ketonic 10mg Tablet", "20's")
    16 = @("Ketonic 30mg Injection", "5 's")
    18 = @("Kynol TR 100mg Capsule", "50 's")
    19 = @("Kynol TR 200mg Capsule", "30 's")
    24 = @("Zithrox 30ml Dry Suspension", "30ml")
    25 = @("Zithrox 15ml Suspension", "15 ml")
    26 = @("Zithrox 250mg Tablet - 6's", "6's")
    27 = @("Zithrox 500mg Tablet", "6 's")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}
